# Errata update: corrected RCP8.5 (with CO2) Iron Availability figures
# for the "2050climate" upper scenario block (columns G:I, rows 15-23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G15" = 27.465440241099206
    "H15" = 26.822806633942893
    "I15" = 27.891794959744313
    "G16" = 31.616444051595622
    "H16" = 30.845295670417421
    "I16" = 32.098159303059035
    "G17" = 23.580524686382418
    "H17" = 23.173799101021263
    "I17" = 23.865225650771166
    "G18" = 26.201010597610512
    "H18" = 25.76995879024981
    "I18" = 26.518647181471671
    "G19" = 22.346410020497366
    "H19" = 21.786779655488136
    "I19" = 22.7575084107545
    "G20" = 27.856106769563876
    "H20" = 27.293586704791668
    "I20" = 28.284540620987599
    "G21" = 25.068820152899647
    "H21" = 24.658640318375713
    "I21" = 25.355659639788648
    "G22" = 25.674992310269783
    "H22" = 24.93134654819676
    "I22" = 26.12625215623542
    "G23" = 31.887886867539834
    "H23" = 31.1180785541275
    "I23" = 32.389130696291296
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Reflect the editor's final selection after making the corrections.
$ws.Range("G15:I23").Select()
